$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 900242908
$ws.Range("B8").Value = "conjunto residencial porvenir reservado casas 2"

$ws.Range("D9").Select()
